# Full run for ZEV Jan R2-4: change roboticLibraryPrep (column L) from the
# text "no" to the boolean FALSE for every data row, and update the
# selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column L holds "roboticLibraryPrep" for rows 2-41 (data rows). Convert
# the literal string "no" into an actual boolean FALSE value, and apply a
# dedicated font (Arial 10, black) to the refreshed cells.
for ($r = 2; $r -le 41; $r++) {
    $cell = $ws.Cells.Item($r, 12)
    $cell.Value = $false
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.Color = 0
}

# Update the active selection to match the new state.
$ws.Range("N33").Select()
